# Daily attendance processing - 2026-01-02 10:00:57
# Normalizes the "Recorded By" (column G) values on the
# "Session Analysis Results" sheet: for the rows below, the first and
# last comma-separated entries are swapped (e.g. "System, x" becomes
# "x, System"), matching the upstream attendance-export reordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = "system, backup@backdoor.com, System"
    5   = "backup@backdoor.com, System"
    7   = "admin@admin.com, System"
    8   = "backup@backdoor.com, System"
    11  = "dnasr281@gmail.com, System"
    17  = "dnasr281@gmail.com, System"
    28  = "system, backup@backdoor.com, System"
    31  = "backup@backdoor.com, System"
    33  = "admin@admin.com, System"
    34  = "backup@backdoor.com, System"
    37  = "dnasr281@gmail.com, System"
    43  = "dnasr281@gmail.com, System"
    54  = "system, backup@backdoor.com, System"
    57  = "backup@backdoor.com, System"
    59  = "admin@admin.com, System"
    60  = "backup@backdoor.com, System"
    63  = "dnasr281@gmail.com, System"
    69  = "dnasr281@gmail.com, System"
    80  = "backup@backdoor.com, System"
    81  = "backup@backdoor.com, System"
    82  = "backup@backdoor.com, System"
    93  = "dnasr281@gmail.com, System"
    94  = "dnasr281@gmail.com, System"
    96  = "dnasr281@gmail.com, System"
    106 = "backup@backdoor.com, System"
    107 = "backup@backdoor.com, System"
    108 = "backup@backdoor.com, System"
    119 = "dnasr281@gmail.com, System"
    120 = "dnasr281@gmail.com, System"
    122 = "dnasr281@gmail.com, System"
    132 = "backup@backdoor.com, System"
    133 = "backup@backdoor.com, System"
    134 = "backup@backdoor.com, System"
    145 = "dnasr281@gmail.com, System"
    146 = "dnasr281@gmail.com, System"
    148 = "dnasr281@gmail.com, System"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
